$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "Repayment schedule" sheet: insert a new (blank) column before the
# old "Late" column (column N), shifting Late / heading / Outstanding
# one column to the right (N->O, O->P, P->Q), for "Variable Instalments".
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert()

# Give the newly inserted column a width roughly matching its neighbour.
$ws.Columns.Item(14).ColumnWidth = 9.83

# Select the sheet + a cell in the new layout, and make it the active tab.
# (This also clears the "active tab" flag on whichever sheet -
# "Transactions" - used to carry it.)
$ws.Activate()
$ws.Range("S7").Select()

Write-Output "done"
